# Auto-generated edit script: apply scheduled-runner price/profit refresh
# to the Leve profit-tracking sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 774.7692
$ws.Range("J43").Value = 746.8333
$ws.Range("L43").Value = 746.8333
$ws.Range("N43").Value = -884.8333

$ws.Range("H51").Value = 2665.7144
$ws.Range("I51").Value = 5000
$ws.Range("J51").Value = 1732
$ws.Range("K51").Value = 5000
$ws.Range("L51").Value = 1732
$ws.Range("M51").Value = -4516
$ws.Range("N51").Value = -2700

$ws.Range("H62").Value = 4234
$ws.Range("I62").Value = 3568.647
$ws.Range("J62").Value = 5262.273
$ws.Range("K62").Value = 3568.647
$ws.Range("L62").Value = 5262.273
$ws.Range("M62").Value = -2944.647
$ws.Range("N62").Value = -6510.273

$ws.Range("H65").Value = 4234
$ws.Range("I65").Value = 3568.647
$ws.Range("J65").Value = 5262.273
$ws.Range("K65").Value = 17843.235
$ws.Range("L65").Value = 26311.365
$ws.Range("M65").Value = -14723.235
$ws.Range("N65").Value = -32551.365

$ws.Range("H86").Value = 7788
$ws.Range("I86").Value = 766.6667
$ws.Range("J86").Value = 20426.4
$ws.Range("K86").Value = 766.6667
$ws.Range("L86").Value = 20426.4
$ws.Range("M86").Value = 356.3333
$ws.Range("N86").Value = -22672.4

$ws.Range("H89").Value = 7788
$ws.Range("I89").Value = 766.6667
$ws.Range("J89").Value = 20426.4
$ws.Range("K89").Value = 3833.3335
$ws.Range("L89").Value = 102132
$ws.Range("M89").Value = 1782.6665
$ws.Range("N89").Value = -113364

$ws.Range("H98").Value = 728.3333
$ws.Range("I98").Value = 476.36365
$ws.Range("K98").Value = 476.36365
$ws.Range("M98").Value = 1021.63635

$ws.Range("H122").Value = 728.3333
$ws.Range("I122").Value = 476.36365
$ws.Range("K122").Value = 1429.09095
$ws.Range("M122").Value = 1020.90905

$ws.Range("H132").Value = 3117.1333
$ws.Range("I132").Value = 3687.7827
$ws.Range("J132").Value = 1242.1428
$ws.Range("K132").Value = 11063.3481
$ws.Range("L132").Value = 3726.4284
$ws.Range("M132").Value = -8533.348100000001
$ws.Range("N132").Value = -8786.428400000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1251.091
$ws.Range("I2").Value = 1326.9333
$ws.Range("J2").Value = 1088.5714
$ws.Range("K2").Value = 1326.9333
$ws.Range("L2").Value = 1088.5714
$ws.Range("M2").Value = -1213.9333
$ws.Range("N2").Value = -1314.5714

$ws.Range("H32").Value = 3313.53
$ws.Range("I32").Value = 2906.7236
$ws.Range("K32").Value = 2906.7236
$ws.Range("M32").Value = -2619.7236

$ws.Range("H61").Value = 451297.1
$ws.Range("J61").Value = 6000
$ws.Range("L61").Value = 6000
$ws.Range("N61").Value = -6424

$ws.Range("H74").Value = 43480896
$ws.Range("I74").Value = 47621650
$ws.Range("J74").Value = 2999.5
$ws.Range("K74").Value = 47621650
$ws.Range("L74").Value = 2999.5
$ws.Range("M74").Value = -47620776
$ws.Range("N74").Value = -4747.5

$ws.Range("H77").Value = 43480896
$ws.Range("I77").Value = 47621650
$ws.Range("J77").Value = 2999.5
$ws.Range("K77").Value = 238108250
$ws.Range("L77").Value = 14997.5
$ws.Range("M77").Value = -238103882
$ws.Range("N77").Value = -23733.5

$ws.Range("H102").Value = 1279.5
$ws.Range("I102").Value = 1072.75
$ws.Range("J102").Value = 1899.75
$ws.Range("K102").Value = 1072.75
$ws.Range("L102").Value = 1899.75
$ws.Range("M102").Value = 549.25
$ws.Range("N102").Value = -5143.75

$ws.Range("H116").Value = 1251.091
$ws.Range("I116").Value = 1326.9333
$ws.Range("J116").Value = 1088.5714
$ws.Range("K116").Value = 1326.9333
$ws.Range("L116").Value = 1088.5714
$ws.Range("M116").Value = 967.0667000000001
$ws.Range("N116").Value = -5676.5714

$ws.Range("H132").Value = 13061.5
$ws.Range("I132").Value = 1564.1111
$ws.Range("J132").Value = 64799.75
$ws.Range("K132").Value = 4692.3333
$ws.Range("L132").Value = 194399.25
$ws.Range("M132").Value = -2162.3333
$ws.Range("N132").Value = -199459.25

$ws.Range("H136").Value = 451297.1
$ws.Range("J136").Value = 6000
$ws.Range("L136").Value = 18000
$ws.Range("N136").Value = -23100

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1251.091
$ws.Range("I3").Value = 1326.9333
$ws.Range("J3").Value = 1088.5714
$ws.Range("K3").Value = 1326.9333
$ws.Range("L3").Value = 1088.5714
$ws.Range("M3").Value = -1212.9333
$ws.Range("N3").Value = -1316.5714

$ws.Range("H86").Value = 1600.2439
$ws.Range("I86").Value = 1382.8695
$ws.Range("J86").Value = 1878
$ws.Range("K86").Value = 1382.8695
$ws.Range("L86").Value = 1878
$ws.Range("M86").Value = -259.8695
$ws.Range("N86").Value = -4124

$ws.Range("H89").Value = 1600.2439
$ws.Range("I89").Value = 1382.8695
$ws.Range("J89").Value = 1878
$ws.Range("K89").Value = 6914.3475
$ws.Range("L89").Value = 9390
$ws.Range("M89").Value = -1298.3475
$ws.Range("N89").Value = -20622

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 367.1
$ws.Range("I22").Value = 271.125
$ws.Range("J22").Value = 751
$ws.Range("K22").Value = 271.125
$ws.Range("L22").Value = 751
$ws.Range("M22").Value = 78.875
$ws.Range("N22").Value = -1451

$ws.Range("H132").Value = 2320.875
$ws.Range("J132").Value = 11004
$ws.Range("L132").Value = 33012
$ws.Range("N132").Value = -38072

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H38").Value = 71428690
$ws.Range("I38").Value = 72.5
$ws.Range("J38").Value = 100000130
$ws.Range("K38").Value = 217.5
$ws.Range("L38").Value = 300000390
$ws.Range("M38").Value = 129.5
$ws.Range("N38").Value = -300001084

$ws.Range("H40").Value = 155
$ws.Range("I40").Value = 82.5
$ws.Range("K40").Value = 330
$ws.Range("M40").Value = -261

$ws.Range("H64").Value = 5309
$ws.Range("I64").Value = 2906
$ws.Range("J64").Value = 6510.5
$ws.Range("K64").Value = 8718
$ws.Range("L64").Value = 19531.5
$ws.Range("M64").Value = -8448
$ws.Range("N64").Value = -20071.5

$ws.Range("H67").Value = 5309
$ws.Range("I67").Value = 2906
$ws.Range("J67").Value = 6510.5
$ws.Range("K67").Value = 8718
$ws.Range("L67").Value = 19531.5
$ws.Range("M67").Value = -7782
$ws.Range("N67").Value = -21403.5

$ws.Range("H107").Value = 4975.2383
$ws.Range("J107").Value = 241.81818
$ws.Range("L107").Value = 725.4545400000001
$ws.Range("N107").Value = -4565.45454

$ws.Range("H131").Value = 257226.56
$ws.Range("J131").Value = 313346.12
$ws.Range("L131").Value = 940038.36
$ws.Range("N131").Value = -950118.36

$ws.Range("H136").Value = 1903.4445
$ws.Range("I136").Value = 1297.1428
$ws.Range("J136").Value = 4025.5
$ws.Range("K136").Value = 3891.4284
$ws.Range("L136").Value = 12076.5
$ws.Range("M136").Value = 1208.5716
$ws.Range("N136").Value = -22276.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 102566680
$ws.Range("I122").Value = 41668770
$ws.Range("J122").Value = 200003340
$ws.Range("K122").Value = 125006310
$ws.Range("L122").Value = 600010020
$ws.Range("M122").Value = -125003860
$ws.Range("N122").Value = -600014920

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2846.75
$ws.Range("I22").Value = 2220.4
$ws.Range("J22").Value = 3294.1428
$ws.Range("K22").Value = 2220.4
$ws.Range("L22").Value = 3294.1428
$ws.Range("M22").Value = -1925.4
$ws.Range("N22").Value = -3884.1428

$ws.Range("H27").Value = 2846.75
$ws.Range("I27").Value = 2220.4
$ws.Range("J27").Value = 3294.1428
$ws.Range("K27").Value = 2220.4
$ws.Range("L27").Value = 3294.1428
$ws.Range("M27").Value = -2113.4
$ws.Range("N27").Value = -3508.1428

$ws.Range("H40").Value = 4407.6665
$ws.Range("I40").Value = 3913.9443
$ws.Range("K40").Value = 3913.9443
$ws.Range("M40").Value = -3777.9443

$ws.Range("H55").Value = 191.72
$ws.Range("I55").Value = 145.08333
$ws.Range("K55").Value = 145.08333
$ws.Range("M55").Value = 27.91667000000001

$ws.Range("H122").Value = 855334.8
$ws.Range("I122").Value = 1403064.4
$ws.Range("J122").Value = 3311
$ws.Range("K122").Value = 4209193.199999999
$ws.Range("L122").Value = 9933
$ws.Range("M122").Value = -4206743.199999999
$ws.Range("N122").Value = -14833

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H30").Value = 0
$ws.Range("J30").Value = 0
$ws.Range("L30").Value = 0
$ws.Range("N30").ClearContents()

$ws.Range("H81").Value = 62501210
$ws.Range("I81").Value = 1160
$ws.Range("J81").Value = 250001340
$ws.Range("K81").Value = 2320
$ws.Range("L81").Value = 500002680
$ws.Range("M81").Value = -1259
$ws.Range("N81").Value = -500004802

$ws.Range("H84").Value = 62501210
$ws.Range("I84").Value = 1160
$ws.Range("J84").Value = 250001340
$ws.Range("K84").Value = 11600
$ws.Range("L84").Value = 2500013400
$ws.Range("M84").Value = -6296
$ws.Range("N84").Value = -2500024008

$ws.Range("H136").Value = 21741704
$ws.Range("I136").Value = 31251192
$ws.Range("J136").Value = 5729.2856
$ws.Range("K136").Value = 93753576
$ws.Range("L136").Value = 17187.8568
$ws.Range("M136").Value = -93751026
$ws.Range("N136").Value = -22287.8568
